# Apply updated "想去人数" (column F) and "最低票价" (column G) values
# to the 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4) worksheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$wsExhibit.Range("F4").Value  = 10364
$wsExhibit.Range("G7").Value  = 150
$wsExhibit.Range("F11").Value = 282
$wsExhibit.Range("F12").Value = 6
$wsExhibit.Range("F13").Value = 12960
$wsExhibit.Range("F18").Value = 155
$wsExhibit.Range("F19").Value = 197
$wsExhibit.Range("F20").Value = 2781
$wsExhibit.Range("F21").Value = 51
$wsExhibit.Range("F23").Value = 2120
$wsExhibit.Range("F24").Value = 142
$wsExhibit.Range("F25").Value = 22
$wsExhibit.Range("F28").Value = 2180
$wsExhibit.Range("F29").Value = 1169
$wsExhibit.Range("F30").Value = 4363
$wsExhibit.Range("F32").Value = 3954
$wsExhibit.Range("F33").Value = 983
$wsExhibit.Range("F34").Value = 2693
$wsExhibit.Range("F35").Value = 3119
$wsExhibit.Range("F37").Value = 1419
$wsExhibit.Range("F38").Value = 222
$wsExhibit.Range("F40").Value = 63
$wsExhibit.Range("F41").Value = 174
$wsExhibit.Range("F42").Value = 634
$wsExhibit.Range("F43").Value = 913
$wsExhibit.Range("F46").Value = 381
$wsExhibit.Range("F47").Value = 131
$wsExhibit.Range("F48").Value = 208
$wsExhibit.Range("F49").Value = 243

# --- 演出 (sheet2) ---
$wsShow.Range("F5").Value  = 73
$wsShow.Range("F9").Value  = 45
$wsShow.Range("F16").Value = 38
$wsShow.Range("F19").Value = 42

# --- 全部类型 (sheet4) ---
$wsAll.Range("F4").Value  = 10364
$wsAll.Range("G6").Value  = 150
$wsAll.Range("F10").Value = 282
$wsAll.Range("F11").Value = 12960
$wsAll.Range("F15").Value = 73
$wsAll.Range("F18").Value = 155
$wsAll.Range("F19").Value = 197
$wsAll.Range("F20").Value = 2781
$wsAll.Range("F21").Value = 2120
$wsAll.Range("F22").Value = 142
$wsAll.Range("F23").Value = 22
$wsAll.Range("F26").Value = 2180
$wsAll.Range("F27").Value = 1169
$wsAll.Range("F30").Value = 4363
$wsAll.Range("F31").Value = 3954
$wsAll.Range("F32").Value = 983
$wsAll.Range("F33").Value = 2693
$wsAll.Range("F34").Value = 3119
$wsAll.Range("F36").Value = 1419
$wsAll.Range("F37").Value = 222
$wsAll.Range("F39").Value = 63
$wsAll.Range("F40").Value = 174
$wsAll.Range("F41").Value = 634
$wsAll.Range("F42").Value = 42
$wsAll.Range("F43").Value = 913
$wsAll.Range("F46").Value = 381
$wsAll.Range("F47").Value = 131
$wsAll.Range("F48").Value = 208
$wsAll.Range("F49").Value = 243
